$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

$ws.Range("D2").Value = 11.5
$ws.Range("D3").Value = 10.25
$ws.Range("B4").Value = 8.5
$ws.Range("C4").Value = 9.75
$ws.Range("F5").Value = 10.2
$ws.Range("H5").Value = 8.5
$ws.Range("E6").Value = 9.800000000000001
$ws.Range("G6").Value = 10.33
$ws.Range("F7").Value = 9.67
$ws.Range("H7").Value = 10.07
$ws.Range("I7").Value = 8
$ws.Range("E8").Value = 11.5
$ws.Range("G8").Value = 9.93
$ws.Range("G9").Value = 12
